$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "sub-label" rows for each category so they repeat the category name,
# e.g. "     New nominations" -> "     Civilian (other than lists), New nominations"
$ws.Range("A7").Value  = "     Civilian (other than lists), New nominations"
$ws.Range("A8").Value  = "     Civilian (other than lists), Confirmed "
$ws.Range("A9").Value  = "     Civilian (other than lists), Unconfirmed "
$ws.Range("A10").Value = "     Civilian (other than lists), Withdrawn "

$ws.Range("A12").Value = "     Civilian (FS, PHS, CG, NOAA), New nominations"
$ws.Range("A13").Value = "     Civilian (FS, PHS, CG, NOAA), Confirmed "
$ws.Range("A14").Value = "     Civilian (FS, PHS, CG, NOAA), Unconfirmed "

$ws.Range("A16").Value = "     Air Force, New nominations"
$ws.Range("A17").Value = "     Air Force, Confirmed "
$ws.Range("A18").Value = "     Air Force, Unconfirmed "

$ws.Range("A20").Value = "     Army, New nominations"
$ws.Range("A21").Value = "     Army, Confirmed "
$ws.Range("A22").Value = "     Army, Unconfirmed "

$ws.Range("A24").Value = "     Navy, New nominations"
$ws.Range("A25").Value = "     Navy, Confirmed "
$ws.Range("A26").Value = "     Navy, Unconfirmed "

$ws.Range("A28").Value = "     Marine Corps, New nominations"
$ws.Range("A29").Value = "     Marine Corps, Confirmed "
$ws.Range("A30").Value = "     Marine Corps, Unconfirmed "

# Remove the bare "Summary" header row (old row 31); the totals below it shift up
# to take its place, preserving their existing values/number formats.
$ws.Rows.Item(31).Delete()

# The (now shifted-up) totals rows need a couple of label tweaks.
$ws.Range("A31").Value = "Total new nominations"
$ws.Range("A35").Value = "Total returned"
